# Update market-price / profit figures on several Leve sheets (ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market board data
# pulled by the scheduled runner. Only the price/profit columns
# (H:N - currentAveragePrice.. LeveProfitHQ) change; some rows gain or
# lose a value in columns M/N depending on whether that row currently
# nets a profit on NQ/HQ crafts.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 250
$ws.Range("I2").Value = 250
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 250
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -137
$ws.Range("N2").ClearContents()

# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()

# Row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()

# Row 29
$ws.Range("H29").Value = 3857.5264
$ws.Range("I29").Value = 3487.8235
$ws.Range("J29").Value = 7000
$ws.Range("K29").Value = 10463.4705
$ws.Range("L29").Value = 21000
$ws.Range("M29").Value = -10182.4705

# Row 31
$ws.Range("H31").Value = 421
$ws.Range("I31").Value = 131.5
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 394.5
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -164.5
$ws.Range("N31").Value = -3460

# Row 41
$ws.Range("H41").Value = 1370.091
$ws.Range("I41").Value = 331.57144
$ws.Range("J41").Value = 3187.5
$ws.Range("K41").Value = 331.57144
$ws.Range("L41").Value = 3187.5
$ws.Range("M41").Value = 108.42856
$ws.Range("N41").Value = -4067.5

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()

# Row 100
$ws.Range("H100").Value = 1734.9
$ws.Range("I100").Value = 1761
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 1761
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -1220
$ws.Range("N100").Value = -2582

# Row 132
$ws.Range("H132").Value = 14728.5
$ws.Range("I132").Value = 13379.375
$ws.Range("J132").Value = 20125
$ws.Range("K132").Value = 40138.125
$ws.Range("L132").Value = 60375
$ws.Range("M132").Value = -37608.125
$ws.Range("N132").Value = -65435

$ws = $wb.Worksheets.Item("ARM")
# Row 21
$ws.Range("H21").Value = 7946.25
$ws.Range("I21").Value = 8928.333000000001
$ws.Range("J21").Value = 5000
$ws.Range("K21").Value = 8928.333000000001
$ws.Range("L21").Value = 5000
$ws.Range("M21").Value = -8554.333000000001
$ws.Range("N21").Value = -5748

# Row 97
$ws.Range("H97").Value = 703.75
$ws.Range("I97").Value = 439.2143
$ws.Range("J97").Value = 2555.5
$ws.Range("K97").Value = 439.2143
$ws.Range("L97").Value = 2555.5
$ws.Range("M97").Value = 56.78570000000002
$ws.Range("N97").Value = -3547.5

# Row 141
$ws.Range("H141").Value = 125999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 125999
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 125999
$ws.Range("N141").Value = -136359

$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 19611.666
$ws.Range("I26").Value = 19611.666
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 19611.666
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -19319.666

# Row 36
$ws.Range("H36").Value = 721.5
$ws.Range("I36").Value = 721.5
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 721.5
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -187.5
$ws.Range("N36").ClearContents()

# Row 135
$ws.Range("H135").Value = 1084544.2
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1084544.2
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 1084544.2
$ws.Range("N135").Value = -1094684.2

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 3750
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3750
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3750
$ws.Range("N22").Value = -4450
$ws.Range("M22").ClearContents()

# Row 39
$ws.Range("H39").Value = 500
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 500
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -109

# Row 44
$ws.Range("H44").Value = 7000
$ws.Range("I44").Value = 7000
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 7000
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -6558

# Row 49
$ws.Range("H49").Value = 500
$ws.Range("I49").Value = 500
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 500
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -318

# Row 68
$ws.Range("H68").Value = 54167.25
$ws.Range("I68").Value = 25000
$ws.Range("J68").Value = 63889.668
$ws.Range("K68").Value = 25000
$ws.Range("L68").Value = 63889.668
$ws.Range("M68").Value = -24251
$ws.Range("N68").Value = -65387.668

# Row 71
$ws.Range("H71").Value = 54167.25
$ws.Range("I71").Value = 25000
$ws.Range("J71").Value = 63889.668
$ws.Range("K71").Value = 75000
$ws.Range("L71").Value = 191669.004
$ws.Range("M71").Value = -71256
$ws.Range("N71").Value = -199157.004

# Row 74
$ws.Range("H74").Value = 53755
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 53755
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 53755
$ws.Range("N74").Value = -55503

# Row 77
$ws.Range("H77").Value = 53755
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 53755
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 161265
$ws.Range("N77").Value = -170001

# Row 134
$ws.Range("H134").Value = 2211.8125
$ws.Range("I134").Value = 2211.8125
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6635.4375
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4100.4375
$ws.Range("N134").ClearContents()

# Row 140
$ws.Range("H140").Value = 55780
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 55780
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 55780
$ws.Range("N140").Value = -66140

$ws = $wb.Worksheets.Item("CUL")
# Row 15
$ws.Range("H15").Value = 227.07692
$ws.Range("I15").Value = 272
$ws.Range("J15").Value = 174.66667
$ws.Range("K15").Value = 816
$ws.Range("L15").Value = 524.00001
$ws.Range("M15").Value = -676
$ws.Range("N15").Value = -804.00001

# Row 50
$ws.Range("H50").Value = 13438
$ws.Range("I50").Value = 14786.286
$ws.Range("J50").Value = 4000
$ws.Range("K50").Value = 44358.858
$ws.Range("L50").Value = 12000
$ws.Range("M50").Value = -43877.858
$ws.Range("N50").Value = -12962

# Row 53
$ws.Range("H53").Value = 13438
$ws.Range("I53").Value = 14786.286
$ws.Range("J53").Value = 4000
$ws.Range("K53").Value = 44358.858
$ws.Range("L53").Value = 12000
$ws.Range("M53").Value = -43877.858
$ws.Range("N53").Value = -12962

# Row 61
$ws.Range("H61").Value = 195.625
$ws.Range("I61").Value = 195.625
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 586.875
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -371.875

# Row 80
$ws.Range("H80").Value = 4127.6665
$ws.Range("I80").Value = 3999.7222
$ws.Range("J80").Value = 4383.5557
$ws.Range("K80").Value = 11999.1666
$ws.Range("L80").Value = 13150.6671
$ws.Range("M80").Value = -11063.1666
$ws.Range("N80").Value = -15022.6671

# Row 82
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

# Row 83
$ws.Range("H83").Value = 4127.6665
$ws.Range("I83").Value = 3999.7222
$ws.Range("J83").Value = 4383.5557
$ws.Range("K83").Value = 35997.49980000001
$ws.Range("L83").Value = 39452.0013
$ws.Range("M83").Value = -31317.49980000001
$ws.Range("N83").Value = -48812.0013

# Row 85
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# Row 135
$ws.Range("H135").Value = 94666.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 94666.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 94666.5
$ws.Range("N135").Value = -104806.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 443.5
$ws.Range("I22").Value = 549.5
$ws.Range("J22").Value = 337.5
$ws.Range("K22").Value = 549.5
$ws.Range("L22").Value = 337.5
$ws.Range("M22").Value = -254.5
$ws.Range("N22").Value = -927.5

# Row 27
$ws.Range("H27").Value = 443.5
$ws.Range("I27").Value = 549.5
$ws.Range("J27").Value = 337.5
$ws.Range("K27").Value = 549.5
$ws.Range("L27").Value = 337.5
$ws.Range("M27").Value = -442.5
$ws.Range("N27").Value = -551.5

# Row 76
$ws.Range("H76").Value = 18893.75
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 18893.75
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 18893.75
$ws.Range("N76").Value = -19569.75

# Row 79
$ws.Range("H79").Value = 18893.75
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 18893.75
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 18893.75
$ws.Range("N79").Value = -21233.75

# Row 93
$ws.Range("H93").Value = 995
$ws.Range("I93").Value = 995
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 995
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 253

# Row 122
$ws.Range("H122").Value = 2846.1538
$ws.Range("I122").Value = 2712.2856
$ws.Range("J122").Value = 3002.3333
$ws.Range("K122").Value = 8136.8568
$ws.Range("L122").Value = 9006.999899999999
$ws.Range("M122").Value = -5686.8568
$ws.Range("N122").Value = -13906.9999

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 127
$ws.Range("H127").Value = 55750
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 55750
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 55750
$ws.Range("N127").Value = -65670

# Row 132
$ws.Range("H132").Value = 4615.5
$ws.Range("I132").Value = 4738
$ws.Range("J132").Value = 4003
$ws.Range("K132").Value = 14214
$ws.Range("L132").Value = 12009
$ws.Range("M132").Value = -11684

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 10000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -22122
$ws.Range("M81").ClearContents()

# Row 84
$ws.Range("H84").Value = 10000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 100000
$ws.Range("N84").Value = -110608
$ws.Range("M84").ClearContents()

# Row 98
$ws.Range("H98").Value = 39295
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 39295
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 39295
$ws.Range("N98").Value = -45285

# Row 100
$ws.Range("H100").Value = 850.3077
$ws.Range("I100").Value = 414
$ws.Range("J100").Value = 3250
$ws.Range("K100").Value = 828
$ws.Range("L100").Value = 6500
$ws.Range("M100").Value = -287
$ws.Range("N100").Value = -7582

# Row 125
$ws.Range("H125").Value = 20000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 20000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 20000
$ws.Range("N125").Value = -29840

# Row 132
$ws.Range("H132").Value = 1660.091
$ws.Range("I132").Value = 1490.8334
$ws.Range("J132").Value = 2421.75
$ws.Range("K132").Value = 4472.5002
$ws.Range("L132").Value = 7265.25
$ws.Range("M132").Value = -1942.5002

# Row 133
$ws.Range("H133").Value = 120000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 120000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 120000
$ws.Range("N133").Value = -130120
